$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (date number format) from the last existing data row
# down onto the three new rows, then fill in the new values.
$ws.Range("A54").Copy()
$ws.Range("A55:A57").PasteSpecial(-4122)

$ws.Range("A55").Value = 46044
$ws.Range("B55").Value = 84

$ws.Range("A56").Value = 46043
$ws.Range("B56").Value = 43

$ws.Range("A57").Value = 46042
$ws.Range("B57").Value = 56

# Move the active selection to G50 (single cell) as in the target workbook
$ws.Range("G50").Select()
